# Mexico Liga MX workbook update (30-05-2024 12:21)
#
# The source data for a handful of fixtures had their rows duplicated/offset
# by one position; this swaps each affected pair of rows back into the
# correct order. For every pair below, the full record (every column from
# "id" (B) through "PL_AhUnder" (AD)) is exchanged between the two rows,
# while the leading row-number column (A) is left untouched since it is
# just the sequential record index.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(2, 3),
    @(36, 37),
    @(72, 73),
    @(130, 131),
    @(251, 252),
    @(264, 265),
    @(310, 311),
    @(318, 319),
    @(322, 323)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rangeA = $ws.Range("B$r1`:AD$r1")
    $rangeB = $ws.Range("B$r2`:AD$r2")

    $valuesA = $rangeA.Value()
    $valuesB = $rangeB.Value()

    $rangeA.Value = $valuesB
    $rangeB.Value = $valuesA
}
